# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45175 (2023-09-06) to serial date 45177 (2023-09-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$firstRow = 2
$lastRow = 135
$newValue = 45177

$ws.Range("C$firstRow`:C$lastRow").Value = $newValue
